# repull data, push all data, mean calculation
# Update the dSF (column F) values for several rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = 3
    6  = -1
    8  = -4
    9  = 1
    10 = 6
    11 = 3
    12 = 2
    13 = 2
    14 = -3
    15 = 4
    16 = 1
    18 = 4
    20 = -3
    22 = -1
    23 = 7
    25 = 2
    26 = -3
    28 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
